$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a number but must stay as literal text;
# set NumberFormat to Text before assigning, then restore default style so
# the cell keeps looking like the rest of the sheet (no explicit style index).
$textCells = @(
    "D5", "D6", "D8", "D10", "D12", "D13", "D17", "D19", "D20", "D21", "D22", "D23", "D25", "D26", "D27", "D29", "D31", "D32", "D34", "D35", "D36", "D39", "D43", "D44", "D47", "D48", "D51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '70.104.43'
$ws.Range("E2").Value = '  -2.59%  '
$ws.Range("D3").Value = '2.521.39'
$ws.Range("E3").Value = '  -5.01%  '
$ws.Range("E4").Value = '  -0.18%  '
$ws.Range("D5").Value = '574.08'
$ws.Range("E5").Value = '  -3.86%  '
$ws.Range("D6").Value = '169.46'
$ws.Range("E6").Value = '  -2.70%  '
$ws.Range("E7").Value = '  -0.09%  '
$ws.Range("D8").Value = '0.508'
$ws.Range("E8").Value = '  -2.85%  '
$ws.Range("D9").Value = '2.520.97'
$ws.Range("E9").Value = '  -4.89%  '
$ws.Range("D10").Value = '0.160'
$ws.Range("E10").Value = '  -4.64%  '
$ws.Range("E11").Value = '  -0.86%  '
$ws.Range("D12").Value = '0.342'
$ws.Range("E12").Value = '  -3.93%  '
$ws.Range("D13").Value = '4.79'
$ws.Range("E13").Value = '  -4.27%  '
$ws.Range("D14").Value = '2.990.16'
$ws.Range("D15").Value = '69.998.93'
$ws.Range("E15").Value = '  -2.83%  '
$ws.Range("E16").Value = '  -2.86%  '
$ws.Range("D17").Value = '24.79'
$ws.Range("E17").Value = '  -5.51%  '
$ws.Range("D18").Value = '2.519.13'
$ws.Range("E18").Value = '  -5.18%  '
$ws.Range("D19").Value = '11.41'
$ws.Range("E19").Value = '  -6.56%  '
$ws.Range("D20").Value = '7.52'
$ws.Range("E20").Value = '  -8.66%  '
$ws.Range("D21").Value = '353.36'
$ws.Range("E21").Value = '  -4.86%  '
$ws.Range("D22").Value = '3.92'
$ws.Range("E22").Value = '  -5.83%  '
$ws.Range("D23").Value = '1.98'
$ws.Range("E23").Value = '  -2.49%  '
$ws.Range("D25").Value = '68.81'
$ws.Range("E25").Value = '  -4.55%  '
$ws.Range("D26").Value = '4.05'
$ws.Range("D27").Value = '9.22'
$ws.Range("E27").Value = '  -5.53%  '
$ws.Range("D28").Value = '2.652.76'
$ws.Range("E28").Value = '  -5.23%  '
$ws.Range("D29").Value = '1.02'
$ws.Range("E29").Value = '  +2.04%  '
$ws.Range("D30").Value = '0.0₃0909'
$ws.Range("E30").Value = '  -6.06%  '
$ws.Range("D31").Value = '7.82'
$ws.Range("E31").Value = '  -2.92%  '
$ws.Range("D32").Value = '478.50'
$ws.Range("E32").Value = '  -4.58%  '
$ws.Range("E33").Value = '  -1.39%  '
$ws.Range("D34").Value = '1.76'
$ws.Range("E34").Value = '  -3.73%  '
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  -0.03%  '
$ws.Range("D36").Value = '157.26'
$ws.Range("E36").Value = '  -3.49%  '
$ws.Range("E37").Value = '  +3.31%  '
$ws.Range("E38").Value = '  -0.94%  '
$ws.Range("D39").Value = '18.51'
$ws.Range("E39").Value = '  -5.30%  '
$ws.Range("E40").Value = '  -0.01%  '
$ws.Range("E41").Value = '  -5.92%  '
$ws.Range("E42").Value = '  -7.12%  '
$ws.Range("D43").Value = '0.318'
$ws.Range("E43").Value = '  -4.12%  '
$ws.Range("D44").Value = '4.70'
$ws.Range("E44").Value = '  -5.84%  '
$ws.Range("E45").Value = '  -7.16%  '
$ws.Range("E46").Value = '  -3.14%  '
$ws.Range("D47").Value = '141.57'
$ws.Range("E47").Value = '  -9.35%  '
$ws.Range("D48").Value = '3.52'
$ws.Range("E48").Value = '  -5.65%  '
$ws.Range("E49").Value = '  -6.71%  '
$ws.Range("E50").Value = '  -7.15%  '
$ws.Range("D51").Value = '0.597'
$ws.Range("E51").Value = '  -0.96%  '

foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
